# Split the combined "ISIC 20T21: Chemicals and pharmaceutical products"
# column into two separate columns: "ISIC 20: Chemicals" and
# "ISIC 21: Pharmaceuticals" on the ItICM sheet, and flag the new
# Pharmaceuticals column for the "other industries" row (row 9), matching
# the "chemicals" industry (row 5) staying mapped only to the Chemicals
# column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItICM")

# Insert a new column before the old column M (i.e. right after the
# existing "chemicals/pharma" column K), shifting ISIC 22 onward one
# column to the right.
$ws.Columns("L").Insert()

# Relabel the original combined column as the dedicated Chemicals column,
# and label the freshly inserted column as the new Pharmaceuticals column.
# (Pharmaceuticals is written first so it lands earlier in the shared
# string table, matching the source workbook's string order.)
$ws.Range("L1").Value = "ISIC 21: Pharmaceuticals"
$ws.Range("K1").Value = "ISIC 20: Chemicals"

# The new column inherited formatting from column K when it was inserted;
# clear that before writing the (unhighlighted) 0 values that belong to
# most rows.
$ws.Range("L2:L9").ClearFormats()
$ws.Range("L2:L9").Value = 0

# "other industries" (row 9) should now also be flagged for the new
# Pharmaceuticals column, matching the existing style used for every
# other flagged cell in this sheet (yellow fill).
$ws.Range("L9").Value = 1
$ws.Range("L9").Interior.Color = 65535
